$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "Sheet1" to "InstituteCancerResearch"
$ws.Name = "InstituteCancerResearch"

# Update Jeffrey Bamber's email address cell (row 2, column E) to include the
# second address.
$ws.Range("E2").Value = "jeffrey.bamber@icr.ac.uk; jeff.bamber@icr.ac.uk"

# Update the saved cell selection to K4, matching the author's last position.
$ws.Range("K4").Select()
